$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 'MSG: None

MSG: The decision has been recorded as ''no decision'' regarding the movie for Friday.
'
$ws.Range("C3").Value = 'MSG: None

MSG: The decision to acquire the rights to "Oppenheimer" has been made.
'
$ws.Range("C4").Value = 'MSG: None

MSG: The decision has been made to acquire the rights to both movies.
'
$ws.Range("C5").Value = 'MSG: None

MSG: I have recorded the decision to acquire the rights for "Barbie."
'
$ws.Range("C6").Value = 'MSG: None

MSG: The decision has been made to acquire the rights for "Barbie." 
'
$ws.Range("C7").Value = 'MSG: None

MSG: The decision was made to not select a movie for Friday.
'
$ws.Range("C8").Value = 'MSG: None

MSG: The decision about which movie to show on Friday resulted in no conclusion.
'
$ws.Range("C9").Value = 'MSG: None

MSG: The decision has been recorded as no decision regarding the movie for Friday. If there are any further discussions or decisions to be made, please let me know!
'
$ws.Range("C10").Value = 'MSG: None

MSG: No decision was made regarding the movie to be shown on Friday.
'
$ws.Range("C11").Value = 'MSG: None

MSG: The decision has been recorded to acquire the rights for "Barbie."
'
$ws.Range("C12").Value = 'MSG: None

MSG: The decision-making process ended without selecting a movie for Friday, so no acquisition will occur at this time.
'
$ws.Range("C13").Value = 'MSG: None

MSG: The decision is recorded as "no_decision."
'
$ws.Range("C14").Value = 'MSG: None

MSG: No decision was made regarding the movie to be shown on Friday.
'
$ws.Range("C15").Value = 'MSG: None

MSG: The decision process concluded without a selection for Friday’s movie, indicating that no agreement was reached.
'
$ws.Range("C16").Value = 'MSG: None

MSG: The decision has been made to acquire the rights to "Barbie" for the movie to be shown on Friday.
'
$ws.Range("C17").Value = 'MSG: None

MSG: The decision has been recorded as no movie being selected for Friday.
'
$ws.Range("C18").Value = 'MSG: None

MSG: A decision about which movie to show on Friday was not reached.
'
$ws.Range("C19").Value = 'MSG: None

MSG: The decision has been recorded, and there will be no movie selected for Friday.
'
$ws.Range("C20").Value = 'MSG: None

MSG: The decision has been recorded, indicating that no movie was selected during the meeting.
'
$ws.Range("C21").Value = 'MSG: None

MSG: The decision regarding the movie to be shown on Friday resulted in no consensus, so there will be no selected film for the week.
'
$ws.Range("C22").Value = 'MSG: None

MSG: The decision process has concluded without a clear selection for Friday''s movie. Therefore, no movie will be acquired.
'
$ws.Range("C23").Value = 'MSG: None

MSG: The decision has been recorded, indicating that no movie was selected for Friday.
'
$ws.Range("C24").Value = 'MSG: None

MSG: The decision has been recorded as no agreement was reached regarding which movie to show on Friday.
'
$ws.Range("C25").Value = 'MSG: None

MSG: The decision resulted in no agreement about which movie to show on Friday.
'
$ws.Range("C26").Value = 'MSG: None

MSG: The decision has been made to acquire the rights for "Barbie."
'
$ws.Range("C27").Value = 'MSG: None

MSG: I have successfully recorded the decision to acquire the rights for "Barbie." The movie will be shown on Friday.
'
$ws.Range("C28").Value = 'MSG: None

MSG: The function for no decision has been successfully called.
'
$ws.Range("C29").Value = 'MSG: None

MSG: The decision to acquire the rights for a movie was not made, as no consensus was reached during the discussion.
'
$ws.Range("C30").Value = 'MSG: None

MSG: The decision to select a movie for Friday was not made, resulting in no acquisition for any movie rights at this time.
'
$ws.Range("C31").Value = 'MSG: None

MSG: The decision about which movie to acquire was ultimately not reached.
'
$ws.Range("C32").Value = 'MSG: None

MSG: The decision has been recorded, and no movie will be shown on Friday.
'
$ws.Range("C33").Value = 'MSG: None

MSG: The decision has been recorded to acquire the rights for "Barbie" to be shown on Friday.
'
$ws.Range("C34").Value = 'MSG: None

MSG: The decision about which movie to show on Friday was not reached, so there is no acquisition action to take.
'
$ws.Range("C35").Value = 'MSG: None

MSG: The decision to acquire the rights for both movies has been successfully recorded.
'
$ws.Range("C36").Value = 'MSG: None

MSG: The decision has been made to acquire the rights for "Oppenheimer."
'
$ws.Range("C37").Value = 'MSG: None

MSG: The decision has been recorded as a no decision, meaning no specific movie was chosen to be shown on Friday.
'
$ws.Range("C38").Value = 'MSG: None

MSG: The decision has been recorded as no movie has been selected for Friday.
'
$ws.Range("C39").Value = 'MSG: None

MSG: The decision regarding the movie to show on Friday has not been made.
'
$ws.Range("C40").Value = 'MSG: None

MSG: The decision process has ended without an agreement on a movie, so no further action will be taken regarding acquiring rights for a movie.
'
$ws.Range("C42").Value = 'MSG: None

MSG: The decision has been recorded as no decision regarding the movie for Friday.
'
$ws.Range("C43").Value = 'MSG: None

MSG: The decision has been recorded as no decision was reached regarding the movie for Friday.
'
$ws.Range("C44").Value = 'MSG: None

MSG: The decision has been recorded successfully. The movie "Barbie" will be shown on Friday.
'
$ws.Range("C45").Value = 'MSG: None

MSG: The committee did not reach a decision on which movie to show on Friday.
'
$ws.Range("C46").Value = 'MSG: None

MSG: The decision process has concluded without selecting a movie for Friday.
'
$ws.Range("C47").Value = 'MSG: None

MSG: The decision process did not lead to an agreement on the movie for Friday, so no decision has been made.
'
$ws.Range("C48").Value = 'MSG: None

MSG: No movie was selected in this meeting.
'
$ws.Range("C49").Value = 'MSG: None

MSG: No decision can be made about Friday''s movie.
'
$ws.Range("C50").Value = 'MSG: None

MSG: The decision has been recorded, indicating that no choice of a movie for Friday was made.
'
$ws.Range("C51").Value = 'MSG: None

MSG: The rights for both movies have been acquired successfully.
'
$ws.Range("C52").Value = 'MSG: None

MSG: The decision has been recorded as no movie being selected for Friday.
'
$ws.Range("C53").Value = 'MSG: None

MSG: The decision has been recorded as "no decision" regarding the movie to be shown on Friday.
'
$ws.Range("C54").Value = 'MSG: None

MSG: The decision regarding the movie to be shown on Friday resulted in no conclusion, adhering to the guidelines provided.
'
$ws.Range("C55").Value = 'MSG: None

MSG: The decision has been recorded as no movie being selected for Friday.
'
$ws.Range("C56").Value = 'MSG: None

MSG: I''ve recorded the decision indicating that no movie was selected for Friday.
'
$ws.Range("C57").Value = 'MSG: None

MSG: The decision has been made to acquire the rights for "Barbie."
'
$ws.Range("C58").Value = 'MSG: None

MSG: It seems that the decision-making process concluded without reaching a consensus on which movie to show on Friday. Therefore, I have recorded the outcome as no decision being made.
'
$ws.Range("C59").Value = 'MSG: None

MSG: The decision has been recorded, and no movie will be acquired for Friday as a consensus was not reached.
'
$ws.Range("C60").Value = 'MSG: None

MSG: I have recorded the decision to acquire the rights for "Barbie" to be shown on Friday.
'
$ws.Range("C61").Value = 'MSG: None

MSG: No decision was made regarding the movie to be shown on Friday.
'
$ws.Range("C62").Value = 'MSG: None

MSG: The decision has been recorded as "no decision."
'
$ws.Range("C63").Value = 'MSG: None

MSG: The decision about which movie to show on Friday has not been made.
'
$ws.Range("C64").Value = 'MSG: None

MSG: The conversation concluded without reaching a decision about which movie to show on Friday.
'
$ws.Range("C65").Value = 'MSG: None

MSG: The decision has been recorded to acquire the rights for both movies.
'
$ws.Range("C66").Value = 'MSG: None

MSG: The decision was recorded as "no decision."
'
$ws.Range("C67").Value = 'MSG: None

MSG: The decision has been recorded that no movie was selected for Friday.
'

$ws.Range("D18").Value = 'no_decision, '
